# 2021 Excess Mortality Update
# - Adds 2021 Q2/Q3/Q4 data rows (Year-Quarter x Race/Ethnicity) that were
#   previously missing (only 2021 Q1 existed before).
# - Revises a handful of already-published 2020/2021 figures with updated numbers.
# - Grows the "Table3" Excel table (and its AutoFilter) from A1:D55 to A1:D73
#   to cover the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table/autofilter range to fit the new rows (55 -> 73 rows incl. header)
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D73"))

# Rewrite the full data block (row 2 through row 73) with the final values.
# Rows 2-19   : 2017-19/2020/2021 Q1 (2021 Q1 values refreshed slightly)
# Rows 20-37  : 2017-19/2020/2021 Q2 (2021 Q2 rows are newly added)
# Rows 38-55  : 2017-19/2020/2021 Q3 (2021 Q3 rows are newly added)
# Rows 56-73  : 2017-19/2020/2021 Q4 (2021 Q4 rows are newly added)
$ws.Cells.Item(2, 1).Value = "2017-19 Q1"
$ws.Cells.Item(2, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(2, 3).Value = 180.6
$ws.Cells.Item(2, 4).Value = 9.8
$ws.Cells.Item(3, 1).Value = "2020 Q1"
$ws.Cells.Item(3, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(3, 3).Value = 173.9
$ws.Cells.Item(3, 4).Value = 9.7
$ws.Cells.Item(4, 1).Value = "2021 Q1"
$ws.Cells.Item(4, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(4, 3).Value = 248
$ws.Cells.Item(4, 4).Value = 11.2
$ws.Cells.Item(5, 1).Value = "2017-19 Q1"
$ws.Cells.Item(5, 2).Value = "Asian"
$ws.Cells.Item(5, 3).Value = 111.5
$ws.Cells.Item(5, 4).Value = 1.3
$ws.Cells.Item(6, 1).Value = "2020 Q1"
$ws.Cells.Item(6, 2).Value = "Asian"
$ws.Cells.Item(6, 3).Value = 114.4
$ws.Cells.Item(6, 4).Value = 1.3
$ws.Cells.Item(7, 1).Value = "2021 Q1"
$ws.Cells.Item(7, 2).Value = "Asian"
$ws.Cells.Item(7, 3).Value = 165.4
$ws.Cells.Item(7, 4).Value = 1.5
$ws.Cells.Item(8, 1).Value = "2017-19 Q1"
$ws.Cells.Item(8, 2).Value = "Black"
$ws.Cells.Item(8, 3).Value = 231.1
$ws.Cells.Item(8, 4).Value = 3.2
$ws.Cells.Item(9, 1).Value = "2020 Q1"
$ws.Cells.Item(9, 2).Value = "Black"
$ws.Cells.Item(9, 3).Value = 233.5
$ws.Cells.Item(9, 4).Value = 3.1
$ws.Cells.Item(10, 1).Value = "2021 Q1"
$ws.Cells.Item(10, 2).Value = "Black"
$ws.Cells.Item(10, 3).Value = 303.4
$ws.Cells.Item(10, 4).Value = 3.5
$ws.Cells.Item(11, 1).Value = "2017-19 Q1"
$ws.Cells.Item(11, 2).Value = "Latino"
$ws.Cells.Item(11, 3).Value = 139.9
$ws.Cells.Item(11, 4).Value = 1.2
$ws.Cells.Item(12, 1).Value = "2020 Q1"
$ws.Cells.Item(12, 2).Value = "Latino"
$ws.Cells.Item(12, 3).Value = 138.7
$ws.Cells.Item(12, 4).Value = 1.2
$ws.Cells.Item(13, 1).Value = "2021 Q1"
$ws.Cells.Item(13, 2).Value = "Latino"
$ws.Cells.Item(13, 3).Value = 261.1
$ws.Cells.Item(13, 4).Value = 1.5
$ws.Cells.Item(14, 1).Value = "2017-19 Q1"
$ws.Cells.Item(14, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(14, 3).Value = 201.7
$ws.Cells.Item(14, 4).Value = 12.7
$ws.Cells.Item(15, 1).Value = "2020 Q1"
$ws.Cells.Item(15, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(15, 3).Value = 199.8
$ws.Cells.Item(15, 4).Value = 12
$ws.Cells.Item(16, 1).Value = "2021 Q1"
$ws.Cells.Item(16, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(16, 3).Value = 298.8
$ws.Cells.Item(16, 4).Value = 14.3
$ws.Cells.Item(17, 1).Value = "2017-19 Q1"
$ws.Cells.Item(17, 2).Value = "White"
$ws.Cells.Item(17, 3).Value = 177.2
$ws.Cells.Item(17, 4).Value = 0.9
$ws.Cells.Item(18, 1).Value = "2020 Q1"
$ws.Cells.Item(18, 2).Value = "White"
$ws.Cells.Item(18, 3).Value = 166.4
$ws.Cells.Item(18, 4).Value = 0.8
$ws.Cells.Item(19, 1).Value = "2021 Q1"
$ws.Cells.Item(19, 2).Value = "White"
$ws.Cells.Item(19, 3).Value = 193
$ws.Cells.Item(19, 4).Value = 0.9
$ws.Cells.Item(20, 1).Value = "2017-19 Q2"
$ws.Cells.Item(20, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(20, 3).Value = 160.3
$ws.Cells.Item(20, 4).Value = 9.4
$ws.Cells.Item(21, 1).Value = "2020 Q2"
$ws.Cells.Item(21, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(21, 3).Value = 185.3
$ws.Cells.Item(21, 4).Value = 9.9
$ws.Cells.Item(22, 1).Value = "2021 Q2"
$ws.Cells.Item(22, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(22, 3).Value = 189.7
$ws.Cells.Item(22, 4).Value = 9.8
$ws.Cells.Item(23, 1).Value = "2017-19 Q2"
$ws.Cells.Item(23, 2).Value = "Asian"
$ws.Cells.Item(23, 3).Value = 97.6
$ws.Cells.Item(23, 4).Value = 1.2
$ws.Cells.Item(24, 1).Value = "2020 Q2"
$ws.Cells.Item(24, 2).Value = "Asian"
$ws.Cells.Item(24, 3).Value = 113.3
$ws.Cells.Item(24, 4).Value = 1.3
$ws.Cells.Item(25, 1).Value = "2021 Q2"
$ws.Cells.Item(25, 2).Value = "Asian"
$ws.Cells.Item(25, 3).Value = 100.1
$ws.Cells.Item(25, 4).Value = 1.2
$ws.Cells.Item(26, 1).Value = "2017-19 Q2"
$ws.Cells.Item(26, 2).Value = "Black"
$ws.Cells.Item(26, 3).Value = 210.5
$ws.Cells.Item(26, 4).Value = 3.1
$ws.Cells.Item(27, 1).Value = "2020 Q2"
$ws.Cells.Item(27, 2).Value = "Black"
$ws.Cells.Item(27, 3).Value = 246.9
$ws.Cells.Item(27, 4).Value = 3.2
$ws.Cells.Item(28, 1).Value = "2021 Q2"
$ws.Cells.Item(28, 2).Value = "Black"
$ws.Cells.Item(28, 3).Value = 214
$ws.Cells.Item(28, 4).Value = 3
$ws.Cells.Item(29, 1).Value = "2017-19 Q2"
$ws.Cells.Item(29, 2).Value = "Latino"
$ws.Cells.Item(29, 3).Value = 124.5
$ws.Cells.Item(29, 4).Value = 1.1
$ws.Cells.Item(30, 1).Value = "2020 Q2"
$ws.Cells.Item(30, 2).Value = "Latino"
$ws.Cells.Item(30, 3).Value = 155.7
$ws.Cells.Item(30, 4).Value = 1.2
$ws.Cells.Item(31, 1).Value = "2021 Q2"
$ws.Cells.Item(31, 2).Value = "Latino"
$ws.Cells.Item(31, 3).Value = 130.4
$ws.Cells.Item(31, 4).Value = 1.1
$ws.Cells.Item(32, 1).Value = "2017-19 Q2"
$ws.Cells.Item(32, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(32, 3).Value = 184.9
$ws.Cells.Item(32, 4).Value = 12.1
$ws.Cells.Item(33, 1).Value = "2020 Q2"
$ws.Cells.Item(33, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(33, 3).Value = 203.2
$ws.Cells.Item(33, 4).Value = 12.2
$ws.Cells.Item(34, 1).Value = "2021 Q2"
$ws.Cells.Item(34, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(34, 3).Value = 203.3
$ws.Cells.Item(34, 4).Value = 11.8
$ws.Cells.Item(35, 1).Value = "2017-19 Q2"
$ws.Cells.Item(35, 2).Value = "White"
$ws.Cells.Item(35, 3).Value = 156.3
$ws.Cells.Item(35, 4).Value = 0.8
$ws.Cells.Item(36, 1).Value = "2020 Q2"
$ws.Cells.Item(36, 2).Value = "White"
$ws.Cells.Item(36, 3).Value = 158.6
$ws.Cells.Item(36, 4).Value = 0.8
$ws.Cells.Item(37, 1).Value = "2021 Q2"
$ws.Cells.Item(37, 2).Value = "White"
$ws.Cells.Item(37, 3).Value = 149.7
$ws.Cells.Item(37, 4).Value = 0.8
$ws.Cells.Item(38, 1).Value = "2017-19 Q3"
$ws.Cells.Item(38, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(38, 3).Value = 158.8
$ws.Cells.Item(38, 4).Value = 9.3
$ws.Cells.Item(39, 1).Value = "2020 Q3"
$ws.Cells.Item(39, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(39, 3).Value = 199.5
$ws.Cells.Item(39, 4).Value = 10.2
$ws.Cells.Item(40, 1).Value = "2021 Q3"
$ws.Cells.Item(40, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(40, 3).Value = 225.5
$ws.Cells.Item(40, 4).Value = 10.9
$ws.Cells.Item(41, 1).Value = "2017-19 Q3"
$ws.Cells.Item(41, 2).Value = "Asian"
$ws.Cells.Item(41, 3).Value = 93.5
$ws.Cells.Item(41, 4).Value = 1.2
$ws.Cells.Item(42, 1).Value = "2020 Q3"
$ws.Cells.Item(42, 2).Value = "Asian"
$ws.Cells.Item(42, 3).Value = 114.4
$ws.Cells.Item(42, 4).Value = 1.3
$ws.Cells.Item(43, 1).Value = "2021 Q3"
$ws.Cells.Item(43, 2).Value = "Asian"
$ws.Cells.Item(43, 3).Value = 106.4
$ws.Cells.Item(43, 4).Value = 1.2
$ws.Cells.Item(44, 1).Value = "2017-19 Q3"
$ws.Cells.Item(44, 2).Value = "Black"
$ws.Cells.Item(44, 3).Value = 205.7
$ws.Cells.Item(44, 4).Value = 3
$ws.Cells.Item(45, 1).Value = "2020 Q3"
$ws.Cells.Item(45, 2).Value = "Black"
$ws.Cells.Item(45, 3).Value = 255.3
$ws.Cells.Item(45, 4).Value = 3.3
$ws.Cells.Item(46, 1).Value = "2021 Q3"
$ws.Cells.Item(46, 2).Value = "Black"
$ws.Cells.Item(46, 3).Value = 244.7
$ws.Cells.Item(46, 4).Value = 3.2
$ws.Cells.Item(47, 1).Value = "2017-19 Q3"
$ws.Cells.Item(47, 2).Value = "Latino"
$ws.Cells.Item(47, 3).Value = 120.5
$ws.Cells.Item(47, 4).Value = 1.1
$ws.Cells.Item(48, 1).Value = "2020 Q3"
$ws.Cells.Item(48, 2).Value = "Latino"
$ws.Cells.Item(48, 3).Value = 176.5
$ws.Cells.Item(48, 4).Value = 1.3
$ws.Cells.Item(49, 1).Value = "2021 Q3"
$ws.Cells.Item(49, 2).Value = "Latino"
$ws.Cells.Item(49, 3).Value = 153.7
$ws.Cells.Item(49, 4).Value = 1.2
$ws.Cells.Item(50, 1).Value = "2017-19 Q3"
$ws.Cells.Item(50, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(50, 3).Value = 180.6
$ws.Cells.Item(50, 4).Value = 12
$ws.Cells.Item(51, 1).Value = "2020 Q3"
$ws.Cells.Item(51, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(51, 3).Value = 230.4
$ws.Cells.Item(51, 4).Value = 13
$ws.Cells.Item(52, 1).Value = "2021 Q3"
$ws.Cells.Item(52, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(52, 3).Value = 255
$ws.Cells.Item(52, 4).Value = 13.2
$ws.Cells.Item(53, 1).Value = "2017-19 Q3"
$ws.Cells.Item(53, 2).Value = "White"
$ws.Cells.Item(53, 3).Value = 150.5
$ws.Cells.Item(53, 4).Value = 0.8
$ws.Cells.Item(54, 1).Value = "2020 Q3"
$ws.Cells.Item(54, 2).Value = "White"
$ws.Cells.Item(54, 3).Value = 165.6
$ws.Cells.Item(54, 4).Value = 0.9
$ws.Cells.Item(55, 1).Value = "2021 Q3"
$ws.Cells.Item(55, 2).Value = "White"
$ws.Cells.Item(55, 3).Value = 165
$ws.Cells.Item(55, 4).Value = 0.9
$ws.Cells.Item(56, 1).Value = "2017-19 Q4"
$ws.Cells.Item(56, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(56, 3).Value = 164.3
$ws.Cells.Item(56, 4).Value = 9.4
$ws.Cells.Item(57, 1).Value = "2020 Q4"
$ws.Cells.Item(57, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(57, 3).Value = 225.4
$ws.Cells.Item(57, 4).Value = 10.8
$ws.Cells.Item(58, 1).Value = "2021 Q4"
$ws.Cells.Item(58, 2).Value = "American Indian or Alaska Native"
$ws.Cells.Item(58, 3).Value = 237.4
$ws.Cells.Item(58, 4).Value = 11
$ws.Cells.Item(59, 1).Value = "2017-19 Q4"
$ws.Cells.Item(59, 2).Value = "Asian"
$ws.Cells.Item(59, 3).Value = 103.6
$ws.Cells.Item(59, 4).Value = 1.3
$ws.Cells.Item(60, 1).Value = "2020 Q4"
$ws.Cells.Item(60, 2).Value = "Asian"
$ws.Cells.Item(60, 3).Value = 143
$ws.Cells.Item(60, 4).Value = 1.4
$ws.Cells.Item(61, 1).Value = "2021 Q4"
$ws.Cells.Item(61, 2).Value = "Asian"
$ws.Cells.Item(61, 3).Value = 113.8
$ws.Cells.Item(61, 4).Value = 1.3
$ws.Cells.Item(62, 1).Value = "2017-19 Q4"
$ws.Cells.Item(62, 2).Value = "Black"
$ws.Cells.Item(62, 3).Value = 222
$ws.Cells.Item(62, 4).Value = 3.2
$ws.Cells.Item(63, 1).Value = "2020 Q4"
$ws.Cells.Item(63, 2).Value = "Black"
$ws.Cells.Item(63, 3).Value = 283
$ws.Cells.Item(63, 4).Value = 3.4
$ws.Cells.Item(64, 1).Value = "2021 Q4"
$ws.Cells.Item(64, 2).Value = "Black"
$ws.Cells.Item(64, 3).Value = 245.2
$ws.Cells.Item(64, 4).Value = 3.2
$ws.Cells.Item(65, 1).Value = "2017-19 Q4"
$ws.Cells.Item(65, 2).Value = "Latino"
$ws.Cells.Item(65, 3).Value = 129.9
$ws.Cells.Item(65, 4).Value = 1.2
$ws.Cells.Item(66, 1).Value = "2020 Q4"
$ws.Cells.Item(66, 2).Value = "Latino"
$ws.Cells.Item(66, 3).Value = 214.3
$ws.Cells.Item(66, 4).Value = 1.4
$ws.Cells.Item(67, 1).Value = "2021 Q4"
$ws.Cells.Item(67, 2).Value = "Latino"
$ws.Cells.Item(67, 3).Value = 160.5
$ws.Cells.Item(67, 4).Value = 1.2
$ws.Cells.Item(68, 1).Value = "2017-19 Q4"
$ws.Cells.Item(68, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(68, 3).Value = 197.1
$ws.Cells.Item(68, 4).Value = 12.5
$ws.Cells.Item(69, 1).Value = "2020 Q4"
$ws.Cells.Item(69, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(69, 3).Value = 269.9
$ws.Cells.Item(69, 4).Value = 14
$ws.Cells.Item(70, 1).Value = "2021 Q4"
$ws.Cells.Item(70, 2).Value = "Native Hawaiian and other Pacific Islander"
$ws.Cells.Item(70, 3).Value = 239.6
$ws.Cells.Item(70, 4).Value = 12.9
$ws.Cells.Item(71, 1).Value = "2017-19 Q4"
$ws.Cells.Item(71, 2).Value = "White"
$ws.Cells.Item(71, 3).Value = 162.1
$ws.Cells.Item(71, 4).Value = 0.8
$ws.Cells.Item(72, 1).Value = "2020 Q4"
$ws.Cells.Item(72, 2).Value = "White"
$ws.Cells.Item(72, 3).Value = 186
$ws.Cells.Item(72, 4).Value = 0.9
$ws.Cells.Item(73, 1).Value = "2021 Q4"
$ws.Cells.Item(73, 2).Value = "White"
$ws.Cells.Item(73, 3).Value = 172.3
$ws.Cells.Item(73, 4).Value = 0.9
